$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices: A=1 ... S=19. Columns A, Q, R contain literal text that
# would otherwise be auto-parsed by Excel as a date / percentage, so those
# cells are pre-formatted as Text ('@') before the value is assigned.
$textCols = @(1, 17, 18)

$rows = @()
$rows += ,@('2025-07-19', 'Independiente del Valle', 'Aucas', 2, 1, 1338458, 2, 4, 2, 2, 1, 0, 0, 0, 2, 1, '52%', '48%', 'L')
$rows += ,@('2025-07-19', 'Universidad Catolica', 'Cuniburo', 3, 0, 1338460, 4, 1, 1, 0, 0, 0, 0, 0, 3, 0, '63%', '37%', 'L')
$rows += ,@('2025-07-19', 'Orense SC', 'El Nacional', 1, 1, 1338459, 10, 1, 3, 1, 0, 0, 0, 0, 1, 1, '61%', '39%', 'E')
$rows += ,@('2025-07-20', 'Delfin SC', 'Barcelona SC', 0, 1, 1338457, 4, 3, 4, 2, 0, 0, 0, 0, 0, 1, '55%', '45%', 'V')
$rows += ,@('2025-07-20', 'Tecnico Universitario', 'Macara', 1, 0, 1338461, 4, 4, 2, 4, 0, 0, 0, 0, 1, 0, '34%', '66%', 'L')
$rows += ,@('2025-07-20', 'Emelec', 'Mushuc Runa SC', 1, 0, 1338462, 5, 6, 1, 4, 0, 0, 0, 0, 1, 0, '57%', '43%', 'L')
$rows += ,@('2025-07-20', 'LDU de Quito', 'Deportivo Cuenca', 2, 2, 1338464, 0, 1, 1, 1, 0, 0, 0, 0, 2, 2, '59%', '41%', 'E')
$rows += ,@('2025-07-22', 'Libertad', 'Manta FC', 1, 0, 1338463, 6, 3, 3, 1, 0, 0, 0, 0, 1, 0, '51%', '49%', 'L')

$startRow = 161
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($textCols -contains $c) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowData[$c - 1]
    }
}
